$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 450, shifting the existing rows 450:471 down to 451:472
$ws.Rows.Item(450).Insert()

# Populate the newly inserted row 450 with the new weekly record
$ws.Cells.Item(450, 1).Value = 4
$ws.Cells.Item(450, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(450, 3).Value = "Los Lagos"
$ws.Cells.Item(450, 4).Value = 45041
$ws.Cells.Item(450, 5).Value = 10
$ws.Cells.Item(450, 6).Value = 100114014
$ws.Cells.Item(450, 7).Value = "Betarraga"
$ws.Cells.Item(450, 8).Value = "Sin especificar"
$ws.Cells.Item(450, 9).Value = "Primera"
$ws.Cells.Item(450, 10).Value = 1000
$ws.Cells.Item(450, 11).Value = 1200
$ws.Cells.Item(450, 12).Value = 1200
$ws.Cells.Item(450, 13).Value = 1200
$ws.Cells.Item(450, 14).Value = "$/paquete 5 unidades"
$ws.Cells.Item(450, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(450, 16).Value = 240
$ws.Cells.Item(450, 17).Value = 5
$ws.Cells.Item(450, 18).Value = "Hortaliza"
